# Replace embedded line breaks in the "Packaging" column (column D) with a
# single space so each value becomes one line instead of several.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$targets = @(6, 8, 14, 15, 16, 23, 24, 25)

foreach ($row in $targets) {
    $cell = $ws.Cells.Item($row, 4)
    $value = $cell.Value2
    if ($value -ne $null) {
        $newValue = $value -replace "`r`n", " "
        $newValue = $newValue -replace "`n", " "
        $newValue = $newValue -replace "`r", " "
        $cell.Value = $newValue
    }
}
